$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 264.5
$ws.Cells.Item(6, 9).Value = 185.66667
$ws.Cells.Item(6, 10).Value = 501
$ws.Cells.Item(6, 11).Value = 557.00001
$ws.Cells.Item(6, 12).Value = 1503
$ws.Cells.Item(6, 13).Value = -445.00001
$ws.Cells.Item(6, 14).Value = -1727
$ws.Cells.Item(8, 8).Value = 11.714286
$ws.Cells.Item(8, 9).Value = 11.714286
$ws.Cells.Item(8, 11).Value = 35.142858
$ws.Cells.Item(8, 13).Value = 103.857142
$ws.Cells.Item(15, 8).Value = 1247.4584
$ws.Cells.Item(15, 9).Value = 1247.4584
$ws.Cells.Item(15, 11).Value = 3742.3752
$ws.Cells.Item(15, 13).Value = -3573.3752
$ws.Cells.Item(32, 8).Value = 3340.5518
$ws.Cells.Item(32, 9).Value = 3653
$ws.Cells.Item(32, 10).Value = 3275.4583
$ws.Cells.Item(32, 11).Value = 3653
$ws.Cells.Item(32, 12).Value = 3275.4583
$ws.Cells.Item(32, 13).Value = -3327
$ws.Cells.Item(32, 14).Value = -3927.4583
$ws.Cells.Item(33, 8).Value = 45909650
$ws.Cells.Item(33, 9).Value = 666885.9
$ws.Cells.Item(33, 10).Value = 142858420
$ws.Cells.Item(33, 11).Value = 666885.9
$ws.Cells.Item(33, 12).Value = 142858420
$ws.Cells.Item(33, 13).Value = -666656.9
$ws.Cells.Item(33, 14).Value = -142858878
$ws.Cells.Item(38, 8).Value = 5962.393
$ws.Cells.Item(38, 9).Value = 726.0769
$ws.Cells.Item(38, 10).Value = 10500.533
$ws.Cells.Item(38, 11).Value = 2178.2307
$ws.Cells.Item(38, 12).Value = 31501.599
$ws.Cells.Item(38, 13).Value = -1806.2307
$ws.Cells.Item(38, 14).Value = -32245.599
$ws.Cells.Item(40, 8).Value = 6669045.5
$ws.Cells.Item(40, 9).Value = 2218.8
$ws.Cells.Item(40, 10).Value = 20002700
$ws.Cells.Item(40, 11).Value = 2218.8
$ws.Cells.Item(40, 12).Value = 20002700
$ws.Cells.Item(40, 13).Value = -2043.8
$ws.Cells.Item(40, 14).Value = -20003050
$ws.Cells.Item(47, 8).Value = 100000
$ws.Cells.Item(47, 9).Value = 100000
$ws.Cells.Item(47, 11).Value = 100000
$ws.Cells.Item(47, 13).Value = -99028
$ws.Cells.Item(58, 8).Value = 10181.6
$ws.Cells.Item(58, 9).Value = 372.5
$ws.Cells.Item(58, 11).Value = 1117.5
$ws.Cells.Item(58, 13).Value = -967.5
$ws.Cells.Item(62, 8).Value = 5146.4165
$ws.Cells.Item(62, 9).Value = 4705.1816
$ws.Cells.Item(62, 11).Value = 4705.1816
$ws.Cells.Item(62, 13).Value = -4081.1816
$ws.Cells.Item(65, 8).Value = 5146.4165
$ws.Cells.Item(65, 9).Value = 4705.1816
$ws.Cells.Item(65, 11).Value = 23525.908
$ws.Cells.Item(65, 13).Value = -20405.908
$ws.Cells.Item(87, 8).Value = 71458.336
$ws.Cells.Item(87, 10).Value = 73863.63
$ws.Cells.Item(87, 12).Value = 73863.63
$ws.Cells.Item(87, 14).Value = -76359.63
$ws.Cells.Item(90, 8).Value = 71458.336
$ws.Cells.Item(90, 10).Value = 73863.63
$ws.Cells.Item(90, 12).Value = 221590.89
$ws.Cells.Item(90, 14).Value = -234070.89
$ws.Cells.Item(98, 8).Value = 3352.7036
$ws.Cells.Item(98, 9).Value = 3444.6365
$ws.Cells.Item(98, 11).Value = 3444.6365
$ws.Cells.Item(98, 13).Value = -1946.6365
$ws.Cells.Item(99, 8).Value = 7346.5713
$ws.Cells.Item(99, 9).Value = 485.2
$ws.Cells.Item(99, 10).Value = 24500
$ws.Cells.Item(99, 11).Value = 1455.6
$ws.Cells.Item(99, 12).Value = 73500
$ws.Cells.Item(99, 13).Value = 42.40000000000009
$ws.Cells.Item(99, 14).Value = -76496
$ws.Cells.Item(100, 8).Value = 8337351
$ws.Cells.Item(100, 9).Value = 3319
$ws.Cells.Item(100, 10).Value = 20838400
$ws.Cells.Item(100, 11).Value = 3319
$ws.Cells.Item(100, 12).Value = 20838400
$ws.Cells.Item(100, 13).Value = -2778
$ws.Cells.Item(100, 14).Value = -20839482
$ws.Cells.Item(106, 8).Value = 47826.555
$ws.Cells.Item(106, 9).Value = 50430.5
$ws.Cells.Item(106, 11).Value = 50430.5
$ws.Cells.Item(106, 13).Value = -49799.5
$ws.Cells.Item(107, 8).Value = 1839.4286
$ws.Cells.Item(107, 9).Value = 1724.4166
$ws.Cells.Item(107, 10).Value = 2529.5
$ws.Cells.Item(107, 11).Value = 1724.4166
$ws.Cells.Item(107, 12).Value = 2529.5
$ws.Cells.Item(107, 13).Value = 195.5834
$ws.Cells.Item(107, 14).Value = -6369.5
$ws.Cells.Item(112, 8).Value = 1514.0435
$ws.Cells.Item(112, 10).Value = 1514.0435
$ws.Cells.Item(112, 12).Value = 4542.1305
$ws.Cells.Item(112, 14).Value = -6758.1305
$ws.Cells.Item(113, 8).Value = 5866.887
$ws.Cells.Item(113, 9).Value = 6076.4165
$ws.Cells.Item(113, 11).Value = 6076.4165
$ws.Cells.Item(113, 13).Value = -2822.4165
$ws.Cells.Item(122, 8).Value = 3352.7036
$ws.Cells.Item(122, 9).Value = 3444.6365
$ws.Cells.Item(122, 11).Value = 10333.9095
$ws.Cells.Item(122, 13).Value = -7883.9095
$ws.Cells.Item(132, 8).Value = 2990.75
$ws.Cells.Item(132, 9).Value = 3114.2144
$ws.Cells.Item(132, 10).Value = 2126.5
$ws.Cells.Item(132, 11).Value = 9342.643199999999
$ws.Cells.Item(132, 12).Value = 6379.5
$ws.Cells.Item(132, 13).Value = -6812.643199999999
$ws.Cells.Item(132, 14).Value = -11439.5
$ws.Cells.Item(134, 8).Value = 68318.5
$ws.Cells.Item(134, 10).Value = 68318.5
$ws.Cells.Item(134, 12).Value = 68318.5
$ws.Cells.Item(134, 14).Value = -78458.5
$ws.Cells.Item(135, 8).Value = 1041.9796
$ws.Cells.Item(135, 9).Value = 1041.9796
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 11).Value = 9377.8164
$ws.Cells.Item(135, 12).Value = 0
$ws.Cells.Item(135, 14).Value = -6842.8164
$ws.Cells.Item(138, 8).Value = 142859820
$ws.Cells.Item(138, 9).Value = 333335000
$ws.Cells.Item(138, 10).Value = 3440
$ws.Cells.Item(138, 11).Value = 1000005000
$ws.Cells.Item(138, 12).Value = 10320
$ws.Cells.Item(138, 13).Value = -999999860
$ws.Cells.Item(138, 14).Value = -20600
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(15, 8).Value = 16670.666
$ws.Cells.Item(15, 9).Value = 10006
$ws.Cells.Item(15, 10).Value = 30000
$ws.Cells.Item(15, 11).Value = 10006
$ws.Cells.Item(15, 12).Value = 30000
$ws.Cells.Item(15, 13).Value = -9656
$ws.Cells.Item(15, 14).Value = -30700
$ws.Cells.Item(24, 8).Value = 85854.5
$ws.Cells.Item(24, 10).Value = 85854.5
$ws.Cells.Item(24, 12).Value = 85854.5
$ws.Cells.Item(24, 14).Value = -86602.5
$ws.Cells.Item(32, 8).Value = 2092.3044
$ws.Cells.Item(32, 9).Value = 1318.0938
$ws.Cells.Item(32, 10).Value = 12002.2
$ws.Cells.Item(32, 11).Value = 1318.0938
$ws.Cells.Item(32, 12).Value = 12002.2
$ws.Cells.Item(32, 13).Value = -1031.0938
$ws.Cells.Item(32, 14).Value = -12576.2
$ws.Cells.Item(61, 8).Value = 1348.3684
$ws.Cells.Item(61, 9).Value = 1145.5
$ws.Cells.Item(61, 11).Value = 1145.5
$ws.Cells.Item(61, 13).Value = -933.5
$ws.Cells.Item(74, 8).Value = 2242.2886
$ws.Cells.Item(74, 9).Value = 1746.9535
$ws.Cells.Item(74, 11).Value = 1746.9535
$ws.Cells.Item(74, 13).Value = -872.9535000000001
$ws.Cells.Item(77, 8).Value = 2242.2886
$ws.Cells.Item(77, 9).Value = 1746.9535
$ws.Cells.Item(77, 11).Value = 8734.7675
$ws.Cells.Item(77, 13).Value = -4366.7675
$ws.Cells.Item(100, 8).Value = 85854.5
$ws.Cells.Item(100, 10).Value = 85854.5
$ws.Cells.Item(100, 12).Value = 85854.5
$ws.Cells.Item(100, 14).Value = -88018.5
$ws.Cells.Item(111, 8).Value = 49221.5
$ws.Cells.Item(111, 10).Value = 49221.5
$ws.Cells.Item(111, 12).Value = 49221.5
$ws.Cells.Item(111, 14).Value = -57401.5
$ws.Cells.Item(122, 8).Value = 5218
$ws.Cells.Item(122, 9).Value = 4967.7144
$ws.Cells.Item(122, 11).Value = 14903.1432
$ws.Cells.Item(122, 13).Value = -12453.1432
$ws.Cells.Item(132, 8).Value = 14260.294
$ws.Cells.Item(132, 9).Value = 7799.6787
$ws.Cells.Item(132, 10).Value = 44409.832
$ws.Cells.Item(132, 11).Value = 23399.0361
$ws.Cells.Item(132, 12).Value = 133229.496
$ws.Cells.Item(132, 13).Value = -20869.0361
$ws.Cells.Item(132, 14).Value = -138289.496
$ws.Cells.Item(136, 8).Value = 1348.3684
$ws.Cells.Item(136, 9).Value = 1145.5
$ws.Cells.Item(136, 11).Value = 3436.5
$ws.Cells.Item(136, 13).Value = -886.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2935.4783
$ws.Cells.Item(86, 9).Value = 2590.7334
$ws.Cells.Item(86, 10).Value = 3581.875
$ws.Cells.Item(86, 11).Value = 2590.7334
$ws.Cells.Item(86, 12).Value = 3581.875
$ws.Cells.Item(86, 13).Value = -1467.7334
$ws.Cells.Item(86, 14).Value = -5827.875
$ws.Cells.Item(89, 8).Value = 2935.4783
$ws.Cells.Item(89, 9).Value = 2590.7334
$ws.Cells.Item(89, 10).Value = 3581.875
$ws.Cells.Item(89, 11).Value = 12953.667
$ws.Cells.Item(89, 12).Value = 17909.375
$ws.Cells.Item(89, 13).Value = -7337.667000000001
$ws.Cells.Item(89, 14).Value = -29141.375
$ws.Cells.Item(99, 8).Value = 47620356
$ws.Cells.Item(99, 9).Value = 76923980
$ws.Cells.Item(99, 11).Value = 76923980
$ws.Cells.Item(99, 13).Value = -76922482
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1010.5714
$ws.Cells.Item(22, 9).Value = 371.84616
$ws.Cells.Item(22, 11).Value = 371.84616
$ws.Cells.Item(22, 13).Value = -21.84616
$ws.Cells.Item(26, 8).Value = 10000
$ws.Cells.Item(26, 10).Value = 10000
$ws.Cells.Item(26, 12).Value = 10000
$ws.Cells.Item(26, 14).Value = -10574
$ws.Cells.Item(31, 8).Value = 1709
$ws.Cells.Item(31, 9).Value = 1536.9375
$ws.Cells.Item(31, 10).Value = 2397.25
$ws.Cells.Item(31, 11).Value = 1536.9375
$ws.Cells.Item(31, 12).Value = 2397.25
$ws.Cells.Item(31, 13).Value = -1241.9375
$ws.Cells.Item(31, 14).Value = -2987.25
$ws.Cells.Item(34, 8).Value = 1709
$ws.Cells.Item(34, 9).Value = 1536.9375
$ws.Cells.Item(34, 10).Value = 2397.25
$ws.Cells.Item(34, 11).Value = 1536.9375
$ws.Cells.Item(34, 12).Value = 2397.25
$ws.Cells.Item(34, 13).Value = -1334.9375
$ws.Cells.Item(34, 14).Value = -2801.25
$ws.Cells.Item(94, 8).Value = 4790.875
$ws.Cells.Item(94, 9).Value = 1665
$ws.Cells.Item(94, 10).Value = 10000.667
$ws.Cells.Item(94, 11).Value = 1665
$ws.Cells.Item(94, 12).Value = 10000.667
$ws.Cells.Item(94, 13).Value = -1214
$ws.Cells.Item(94, 14).Value = -10902.667
$ws.Cells.Item(99, 8).Value = 3638.4119
$ws.Cells.Item(99, 9).Value = 1595.125
$ws.Cells.Item(99, 10).Value = 5454.6665
$ws.Cells.Item(99, 11).Value = 1595.125
$ws.Cells.Item(99, 12).Value = 5454.6665
$ws.Cells.Item(99, 13).Value = -97.125
$ws.Cells.Item(99, 14).Value = -8450.666499999999
$ws.Cells.Item(126, 8).Value = 3638.4119
$ws.Cells.Item(126, 9).Value = 1595.125
$ws.Cells.Item(126, 10).Value = 5454.6665
$ws.Cells.Item(126, 11).Value = 4785.375
$ws.Cells.Item(126, 12).Value = 16363.9995
$ws.Cells.Item(126, 13).Value = -2315.375
$ws.Cells.Item(126, 14).Value = -21303.9995
$ws.Cells.Item(132, 8).Value = 2544.6
$ws.Cells.Item(132, 9).Value = 2094
$ws.Cells.Item(132, 11).Value = 6282
$ws.Cells.Item(132, 13).Value = -3752
$ws.Cells.Item(134, 8).Value = 1315.6061
$ws.Cells.Item(134, 9).Value = 1322.4482
$ws.Cells.Item(134, 10).Value = 1266
$ws.Cells.Item(134, 11).Value = 3967.3446
$ws.Cells.Item(134, 12).Value = 3798
$ws.Cells.Item(134, 13).Value = -1432.3446
$ws.Cells.Item(134, 14).Value = -8868
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 240.33333
$ws.Cells.Item(7, 9).Value = 48
$ws.Cells.Item(7, 10).Value = 336.5
$ws.Cells.Item(7, 11).Value = 144
$ws.Cells.Item(7, 12).Value = 1009.5
$ws.Cells.Item(7, 13).Value = -32
$ws.Cells.Item(7, 14).Value = -1233.5
$ws.Cells.Item(37, 8).Value = 106999.11
$ws.Cells.Item(37, 10).Value = 106999.11
$ws.Cells.Item(37, 12).Value = 320997.33
$ws.Cells.Item(37, 14).Value = -321221.33
$ws.Cells.Item(40, 8).Value = 122.14286
$ws.Cells.Item(40, 9).Value = 136.41176
$ws.Cells.Item(40, 10).Value = 61.5
$ws.Cells.Item(40, 11).Value = 545.6470399999999
$ws.Cells.Item(40, 12).Value = 246
$ws.Cells.Item(40, 13).Value = -476.6470399999999
$ws.Cells.Item(40, 14).Value = -384
$ws.Cells.Item(60, 8).Value = 1030.8889
$ws.Cells.Item(60, 9).Value = 44.75
$ws.Cells.Item(60, 11).Value = 134.25
$ws.Cells.Item(60, 13).Value = 116.75
$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 14).Value = 0
$ws.Cells.Item(69, 8).Value = 9999
$ws.Cells.Item(69, 9).Value = 9999
$ws.Cells.Item(69, 11).Value = 29997
$ws.Cells.Item(69, 13).Value = -29186
$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 14).Value = 0
$ws.Cells.Item(72, 8).Value = 9999
$ws.Cells.Item(72, 9).Value = 9999
$ws.Cells.Item(72, 11).Value = 89991
$ws.Cells.Item(72, 13).Value = -85935
$ws.Cells.Item(81, 8).Value = 12566.333
$ws.Cells.Item(81, 9).Value = 10006.5
$ws.Cells.Item(81, 10).Value = 13846.25
$ws.Cells.Item(81, 11).Value = 30019.5
$ws.Cells.Item(81, 12).Value = 41538.75
$ws.Cells.Item(81, 13).Value = -28896.5
$ws.Cells.Item(81, 14).Value = -43784.75
$ws.Cells.Item(84, 8).Value = 12566.333
$ws.Cells.Item(84, 9).Value = 10006.5
$ws.Cells.Item(84, 10).Value = 13846.25
$ws.Cells.Item(84, 11).Value = 90058.5
$ws.Cells.Item(84, 12).Value = 124616.25
$ws.Cells.Item(84, 13).Value = -84442.5
$ws.Cells.Item(84, 14).Value = -135848.25
$ws.Cells.Item(92, 8).Value = 294.07693
$ws.Cells.Item(92, 10).Value = 268
$ws.Cells.Item(92, 12).Value = 804
$ws.Cells.Item(92, 14).Value = -3300
$ws.Cells.Item(122, 8).Value = 1301.16
$ws.Cells.Item(122, 10).Value = 1016.3889
$ws.Cells.Item(122, 12).Value = 9147.500100000001
$ws.Cells.Item(122, 14).Value = -14047.5001
$ws.Cells.Item(137, 8).Value = 2153
$ws.Cells.Item(137, 10).Value = 3446.5
$ws.Cells.Item(137, 12).Value = 10339.5
$ws.Cells.Item(137, 14).Value = -20539.5
$ws.Cells.Item(140, 8).Value = 1320.5
$ws.Cells.Item(140, 9).Value = 911.6667
$ws.Cells.Item(140, 10).Value = 5000
$ws.Cells.Item(140, 11).Value = 2735.0001
$ws.Cells.Item(140, 12).Value = 15000
$ws.Cells.Item(140, 13).Value = 2444.9999
$ws.Cells.Item(140, 14).Value = -25360
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 194.53334
$ws.Cells.Item(2, 9).Value = 202.72728
$ws.Cells.Item(2, 10).Value = 172
$ws.Cells.Item(2, 11).Value = 202.72728
$ws.Cells.Item(2, 12).Value = 172
$ws.Cells.Item(2, 13).Value = -89.72728000000001
$ws.Cells.Item(2, 14).Value = -398
$ws.Cells.Item(17, 8).Value = 4966.6665
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = 4966.6665
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 13).Value = 4966.6665
$ws.Cells.Item(17, 14).Value = -5302.6665
$ws.Cells.Item(19, 8).Value = 20150.666
$ws.Cells.Item(19, 9).Value = 24666
$ws.Cells.Item(19, 10).Value = 15635.333
$ws.Cells.Item(19, 11).Value = 24666
$ws.Cells.Item(19, 12).Value = 15635.333
$ws.Cells.Item(19, 13).Value = -24378
$ws.Cells.Item(19, 14).Value = -16211.333
$ws.Cells.Item(34, 8).Value = 11513.25
$ws.Cells.Item(34, 10).Value = 11513.25
$ws.Cells.Item(34, 12).Value = 11513.25
$ws.Cells.Item(34, 14).Value = -12049.25
$ws.Cells.Item(76, 8).Value = 11513.25
$ws.Cells.Item(76, 10).Value = 11513.25
$ws.Cells.Item(76, 12).Value = 11513.25
$ws.Cells.Item(76, 14).Value = -12143.25
$ws.Cells.Item(79, 8).Value = 11513.25
$ws.Cells.Item(79, 10).Value = 11513.25
$ws.Cells.Item(79, 12).Value = 11513.25
$ws.Cells.Item(79, 14).Value = -13697.25
$ws.Cells.Item(97, 8).Value = 104.5
$ws.Cells.Item(97, 9).Value = 104.5
$ws.Cells.Item(97, 11).Value = 104.5
$ws.Cells.Item(97, 13).Value = 391.5
$ws.Cells.Item(102, 8).Value = 2557.5715
$ws.Cells.Item(102, 9).Value = 1783.1177
$ws.Cells.Item(102, 11).Value = 1783.1177
$ws.Cells.Item(102, 13).Value = -161.1177
$ws.Cells.Item(107, 8).Value = 1543.7333
$ws.Cells.Item(107, 9).Value = 3200
$ws.Cells.Item(107, 10).Value = 715.6
$ws.Cells.Item(107, 11).Value = 3200
$ws.Cells.Item(107, 12).Value = 715.6
$ws.Cells.Item(107, 13).Value = -1280
$ws.Cells.Item(107, 14).Value = -4555.6
$ws.Cells.Item(122, 8).Value = 2398.913
$ws.Cells.Item(122, 10).Value = 3342.2727
$ws.Cells.Item(122, 12).Value = 10026.8181
$ws.Cells.Item(122, 14).Value = -14926.8181
$ws.Cells.Item(132, 8).Value = 3313
$ws.Cells.Item(132, 9).Value = 3306
$ws.Cells.Item(132, 10).Value = 3320
$ws.Cells.Item(132, 11).Value = 9918
$ws.Cells.Item(132, 12).Value = 9960
$ws.Cells.Item(132, 13).Value = -7388
$ws.Cells.Item(132, 14).Value = -15020
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = 0
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 14).Value = 0
$ws.Cells.Item(25, 8).Value = 38503.5
$ws.Cells.Item(25, 9).Value = 38503.5
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 11).Value = 38503.5
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 14).Value = -38273.5
$ws.Cells.Item(33, 8).Value = 74250
$ws.Cells.Item(33, 10).Value = 28500
$ws.Cells.Item(33, 12).Value = 28500
$ws.Cells.Item(33, 14).Value = -29080
$ws.Cells.Item(40, 8).Value = 5510.0527
$ws.Cells.Item(40, 9).Value = 4981.8823
$ws.Cells.Item(40, 11).Value = 4981.8823
$ws.Cells.Item(40, 13).Value = -4845.8823
$ws.Cells.Item(46, 8).Value = 1759.5714
$ws.Cells.Item(46, 9).Value = 1213.4
$ws.Cells.Item(46, 10).Value = 3125
$ws.Cells.Item(46, 11).Value = 1213.4
$ws.Cells.Item(46, 12).Value = 3125
$ws.Cells.Item(46, 13).Value = -1025.4
$ws.Cells.Item(46, 14).Value = -3501
$ws.Cells.Item(59, 8).Value = 73000
$ws.Cells.Item(59, 10).Value = 73000
$ws.Cells.Item(59, 12).Value = 73000
$ws.Cells.Item(59, 14).Value = -74308
$ws.Cells.Item(61, 8).Value = 831.875
$ws.Cells.Item(61, 9).Value = 807.8570999999999
$ws.Cells.Item(61, 11).Value = 807.8570999999999
$ws.Cells.Item(61, 13).Value = -605.8570999999999
$ws.Cells.Item(74, 8).Value = 65000
$ws.Cells.Item(74, 9).Value = 65000
$ws.Cells.Item(74, 11).Value = 65000
$ws.Cells.Item(74, 13).Value = -64002
$ws.Cells.Item(77, 8).Value = 65000
$ws.Cells.Item(77, 9).Value = 65000
$ws.Cells.Item(77, 11).Value = 195000
$ws.Cells.Item(77, 13).Value = -190008
$ws.Cells.Item(93, 8).Value = 2858.4
$ws.Cells.Item(93, 9).Value = 2394.7188
$ws.Cells.Item(93, 10).Value = 4713.125
$ws.Cells.Item(93, 11).Value = 2394.7188
$ws.Cells.Item(93, 12).Value = 4713.125
$ws.Cells.Item(93, 13).Value = -1146.7188
$ws.Cells.Item(93, 14).Value = -7209.125
$ws.Cells.Item(100, 8).Value = 2202.2
$ws.Cells.Item(100, 9).Value = 1214.7778
$ws.Cells.Item(100, 10).Value = 3683.3333
$ws.Cells.Item(100, 11).Value = 1214.7778
$ws.Cells.Item(100, 12).Value = 3683.3333
$ws.Cells.Item(100, 13).Value = -673.7778000000001
$ws.Cells.Item(100, 14).Value = -4765.3333
$ws.Cells.Item(113, 8).Value = 831.875
$ws.Cells.Item(113, 9).Value = 807.8570999999999
$ws.Cells.Item(113, 11).Value = 807.8570999999999
$ws.Cells.Item(113, 13).Value = 1362.1429
$ws.Cells.Item(122, 8).Value = 5168.8237
$ws.Cells.Item(122, 9).Value = 4151
$ws.Cells.Item(122, 10).Value = 6073.5557
$ws.Cells.Item(122, 11).Value = 12453
$ws.Cells.Item(122, 12).Value = 18220.6671
$ws.Cells.Item(122, 13).Value = -10003
$ws.Cells.Item(122, 14).Value = -23120.6671
$ws.Cells.Item(132, 8).Value = 3604.138
$ws.Cells.Item(132, 9).Value = 3304.2126
$ws.Cells.Item(132, 10).Value = 4885.636
$ws.Cells.Item(132, 11).Value = 9912.6378
$ws.Cells.Item(132, 12).Value = 14656.908
$ws.Cells.Item(132, 13).Value = -7382.6378
$ws.Cells.Item(132, 14).Value = -19716.908
$ws.Cells.Item(134, 8).Value = 34184.777
$ws.Cells.Item(134, 10).Value = 34184.777
$ws.Cells.Item(134, 12).Value = 34184.777
$ws.Cells.Item(134, 14).Value = -44324.777
$ws.Cells.Item(135, 8).Value = 0
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 14).Value = 0
$ws.Cells.Item(136, 8).Value = 2466.7
$ws.Cells.Item(136, 9).Value = 2102.0571
$ws.Cells.Item(136, 10).Value = 5019.2
$ws.Cells.Item(136, 11).Value = 6306.1713
$ws.Cells.Item(136, 12).Value = 15057.6
$ws.Cells.Item(136, 13).Value = -3756.1713
$ws.Cells.Item(136, 14).Value = -20157.6
$ws.Cells.Item(137, 8).Value = 0
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 14).Value = 0
$ws.Cells.Item(138, 8).Value = 0
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 14).Value = 0
$ws.Cells.Item(139, 8).Value = 100000
$ws.Cells.Item(139, 10).Value = 100000
$ws.Cells.Item(139, 12).Value = 100000
$ws.Cells.Item(139, 14).Value = -110280
$ws.Cells.Item(140, 8).Value = 58461.54
$ws.Cells.Item(141, 8).Value = 149833.25
$ws.Cells.Item(141, 10).Value = 149833.25
$ws.Cells.Item(141, 12).Value = 149833.25
$ws.Cells.Item(141, 14).Value = -160193.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(11, 8).Value = 62500
$ws.Cells.Item(11, 10).Value = 50000
$ws.Cells.Item(11, 12).Value = 50000
$ws.Cells.Item(11, 14).Value = -50284
$ws.Cells.Item(116, 8).Value = 34500
$ws.Cells.Item(116, 10).Value = 34500
$ws.Cells.Item(116, 12).Value = 34500
$ws.Cells.Item(116, 14).Value = -43678
$ws.Cells.Item(122, 8).Value = 8776105
$ws.Cells.Item(122, 10).Value = 5000.6
$ws.Cells.Item(122, 12).Value = 15001.8
$ws.Cells.Item(122, 14).Value = -19901.8
$ws.Cells.Item(132, 8).Value = 3036.5
$ws.Cells.Item(132, 9).Value = 1929.3334
$ws.Cells.Item(132, 10).Value = 3700.8
$ws.Cells.Item(132, 11).Value = 5788.0002
$ws.Cells.Item(132, 12).Value = 11102.4
$ws.Cells.Item(132, 13).Value = -3258.0002
$ws.Cells.Item(132, 14).Value = -16162.4
$ws.Cells.Item(136, 8).Value = 4440
$ws.Cells.Item(136, 9).Value = 4125.758
$ws.Cells.Item(136, 10).Value = 5736.25
$ws.Cells.Item(136, 11).Value = 12377.274
$ws.Cells.Item(136, 12).Value = 17208.75
$ws.Cells.Item(136, 13).Value = -9827.273999999999
$ws.Cells.Item(136, 14).Value = -22308.75
